$d = $word.ActiveDocument

$replacements = @(
    @{old = "60×71="; new = "80×68="},
    @{old = "49×61="; new = "77×88="},
    @{old = "81×46="; new = "50×79="},
    @{old = "79×48="; new = "75×58="},
    @{old = "37×66="; new = "87×12="},
    @{old = "25×65="; new = "62×80="},
    @{old = "48×41="; new = "91×70="},
    @{old = "88×26="; new = "11×92="},
    @{old = "70×37="; new = "75×75="},
    @{old = "96×90="; new = "21×61="},
    @{old = "28×14="; new = "83×62="},
    @{old = "86×39="; new = "65×42="},
    @{old = "98×43="; new = "84×32="},
    @{old = "16×23="; new = "19×80="},
    @{old = "17×52="; new = "96×60="},
    @{old = "82×29="; new = "37×14="},
    @{old = "57×45="; new = "23×77="},
    @{old = "61×61="; new = "56×37="},
    @{old = "83×90="; new = "24×20="},
    @{old = "11×13="; new = "14×37="},
    @{old = "22×90="; new = "98×51="},
    @{old = "12×63="; new = "19×29="},
    @{old = "20×50="; new = "89×19="},
    @{old = "52×50="; new = "26×81="},
    @{old = "77×17="; new = "67×24="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}

Write-Host "Done applying replacements"
